$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.559.51'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '1.622.27'
$ws.Range("E3").Value = '  -1.03%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.72'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.521'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.22%  '
$ws.Range("E7").Value = '  +0.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.06'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.263'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.07%  '
$ws.Range("E10").Value = '  +0.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0887'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.21%  '
$ws.Range("D12").Value = '1.852.49'
$ws.Range("E12").Value = '  -0.83%  '
$ws.Range("D13").Value = '1.622.00'
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("E14").Value = '  +0.25%  '
$ws.Range("E15").Value = '  -2.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.41'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").Value = '27.565.12'
$ws.Range("E17").Value = '  +0.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '230.34'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.05%  '
$ws.Range("E19").Value = '  -0.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.57'
$ws.Range("D20").ClearFormats()
$ws.Range("E21").Value = '  +0.20%  '
$ws.Range("E22").Value = '  -0.50%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.92'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.07'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +6.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.22'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.88'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.38%  '
$ws.Range("E27").Value = '  -1.36%  '
$ws.Range("E28").Value = '  +0.24%  '
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("E30").Value = '  -0.75%  '
$ws.Range("E31").Value = '  -1.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.28'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.42%  '
$ws.Range("D33").Value = '1.455.62'
$ws.Range("E33").Value = '  +2.70%  '
$ws.Range("E34").Value = '  -3.00%  '
$ws.Range("E35").Value = '  -3.05%  '
$ws.Range("E36").Value = '  -0.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.950'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +6.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.563'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.50%  '
$ws.Range("E39").Value = '  +0.21%  '
$ws.Range("E40").Value = '  -1.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.26'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +6.52%  '
$ws.Range("E42").Value = '  +0.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.01'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.47'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.43'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.22'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.49%  '
$ws.Range("D47").Value = '1.763.56'
$ws.Range("E47").Value = '  -0.74%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.68'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.25'
$ws.Range("D49").ClearFormats()
$ws.Range("E50").Value = '  -0.65%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0986'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.55%  '
